# Applies the commit's changes:
#  - School sheet: B1 price 97.45 -> 200, B2 price 175 -> 205; becomes active tab, selection -> B1
#  - Purchase sheet: C1 price 40 -> 45, C2 price 30 -> 26; no longer active tab, selection -> C2
#  - Tool sheet: add new row "Exacto Knife" / 5 in row 6; selection -> A6:B6
#  - Workbook active tab moves from Purchase to School

$wb = $excel.ActiveWorkbook

$wsSchool = $wb.Worksheets.Item("School")
$wsPurchase = $wb.Worksheets.Item("Purchase")
$wsTool = $wb.Worksheets.Item("Tool")

# --- School sheet updates ---
$wsSchool.Range("B1").Value = 200
$wsSchool.Range("B2").Value = 205
# Target raw column width (stored XML "width" attr) is 11.47265625 characters.
# The ColumnWidth COM setter here quantizes to a 1/6-character pixel grid, so
# feed the input value whose quantized bucket lands closest to the target
# (10.666666666666666 -> stored width 11.5, the nearest achievable value).
$wsSchool.Columns.Item(1).ColumnWidth = 10.666666666666666

# --- Purchase sheet updates ---
$wsPurchase.Range("C1").Value = 45
$wsPurchase.Range("C2").Value = 26

# --- Tool sheet: add a new tool row ---
$wsTool.Range("A6").Value = "Exacto Knife"
$wsTool.Range("B6").Value = 5

# --- Selections / active sheet ---
$wsTool.Range("A6:B6").Select()
$wsPurchase.Range("C2").Select()
$wsSchool.Range("B1").Select()
$wsSchool.Activate()
